# Update deflator comparison figures on the (single) worksheet.
# Sections: current (rows 2-29), difference (rows 30-57), previous (rows 58-85)
# Rows touched here: 12 (Federal Social Benefits Contribution / current),
# 16 (Fiscal Impact / current), 22 (State Contribution / current),
# 29 (State Ui Contribution / current), 40 (Federal Social Benefits
# Contribution / difference), 44 (Fiscal Impact / difference),
# 50 (State Contribution / difference), 57 (State Ui Contribution /
# difference).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = @{ J = 0.0191; K = 0.1269; L = 0.1731; M = -0.1102; N = -0.0957; O = -0.0808; P = -0.066; Q = -0.0309; R = -0.0692 }
    16 = @{ J = -1.488; K = -1.6061; L = -0.5396; M = -0.5624; N = -0.5968; O = -0.1776; P = -0.1055; Q = 0.0429; R = 0.0089 }
    22 = @{ R = 0.2276 }
    29 = @{ J = -0.0275; K = -0.0062; L = 0.0071; N = 0.0122; O = 0.0075; Q = -0.0016; R = -0.004 }
    40 = @{ J = 0.0049; K = 0.003; L = 0.0009; M = -0.003; N = 0.0145; O = 0.0142; P = 0.0139; Q = 0.0136; R = 0.0105 }
    44 = @{ J = 0.0471; K = 0.0335; L = 0.0344; M = 0.0075; N = 0.0215; O = 0.018; P = 0.0343; Q = 0.0175; R = 0.0609 }
    50 = @{ R = 0.0463 }
    57 = @{ J = -0.0007; K = 0.0005; L = 0.0004; M = 0.0004; N = 0.0005; O = 0.0003; R = 0.0002 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
